# Applies the "Trade #20 closed" update described in the commit:
#   - Trade #48 (row index in the trade logs) flips from OPEN -> CLOSED
#     with a small negative P&L (early_exit), which also nudges the
#     Summary / Strategy Status rollups.
#   - A brand-new trade (#81) is appended as OPEN to both the
#     "All Trades" log and the per-strategy "MarketMaking" log.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($cell, [string]$text)
    # Writing a plain .Value for strings that look like dates (e.g.
    # "2026-02-17") makes Excel silently re-interpret them as date
    # serials. Force literal text by flipping the number format to
    # Text for the assignment, then restore the "Normal" style so the
    # cell doesn't carry a lingering explicit style index.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------
# 1) Summary sheet rollups
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1400.4   # Current Capital
$summary.Range("B4").Value = 0.19     # Total P&L $
$summary.Range("B5").Value = 0.08     # Total P&L %
$summary.Range("B6").Value = 48       # Total Trades
$summary.Range("B8").Value = 20       # Losing Trades
$summary.Range("B9").Value = 45.83    # Win Rate %

# ---------------------------------------------------------------
# 2) Strategy Status sheet - MarketMaking row (row 5)
# ---------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 100.4   # Capital
$status.Range("D5").Value = 15      # Trades
$status.Range("E5").Value = 0.08    # P&L $
$status.Range("F5").Value = 0.4     # P&L %
$status.Range("G5").Value = 53.33   # Win Rate %

# ---------------------------------------------------------------
# 3) All Trades sheet
# ---------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Existing trade #48 (sheet row 49) closes out.
$allTrades.Cells.Item(49, 7).Value = 0.171474        # G: Exit Price
$allTrades.Cells.Item(49, 8).Value = "CLOSED"        # H: Status
$allTrades.Cells.Item(49, 9).Value = -14.2629        # I: P&L %
$allTrades.Cells.Item(49, 10).Value = -0.03          # J: P&L $
$allTrades.Cells.Item(49, 11).Value = 100.4          # K: Capital After
$allTrades.Cells.Item(49, 12).Value = "early_exit"   # L: Exit Reason
$allTrades.Cells.Item(49, 13).Value = 0.15           # M: Duration (min)

# New trade #81 appended as sheet row 82.
$allTrades.Cells.Item(82, 1).Value = 81              # A: Trade #
Set-TextValue $allTrades.Cells.Item(82, 2) "2026-02-17"   # B: Date
Set-TextValue $allTrades.Cells.Item(82, 3) "20:53:16"     # C: Time
$allTrades.Cells.Item(82, 4).Value = "MarketMaking"  # D: Strategy
$allTrades.Cells.Item(82, 5).Value = "UP"            # E: Side
$allTrades.Cells.Item(82, 6).Value = 0.2             # F: Entry Price
$allTrades.Cells.Item(82, 8).Value = "OPEN"          # H: Status
$allTrades.Cells.Item(82, 9).Value = 0               # I: P&L %
$allTrades.Cells.Item(82, 10).Value = 0              # J: P&L $
$allTrades.Cells.Item(82, 11).Value = 100.4269627845085 # K: Capital After
$allTrades.Cells.Item(82, 13).Value = 0              # M: Duration (min)
$allTrades.Cells.Item(82, 14).Value = 0              # N: Entry Slippage (bps)
$allTrades.Cells.Item(82, 15).Value = 0              # O: Exit Slippage (bps)
$allTrades.Cells.Item(82, 16).Value = 0.6            # P: Confidence
$allTrades.Cells.Item(82, 17).Value = "Normal spread capture: 19600 bps" # Q: Entry Reason

# ---------------------------------------------------------------
# 4) MarketMaking sheet (per-strategy trade log; column order for
#    L..Q differs from "All Trades")
# ---------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

# Existing trade #48 (sheet row 16) closes out.
$mm.Cells.Item(16, 7).Value = 0.171474       # G: Exit Price
$mm.Cells.Item(16, 8).Value = "CLOSED"       # H: Status
$mm.Cells.Item(16, 9).Value = -14.2629       # I: P&L %
$mm.Cells.Item(16, 10).Value = -0.03         # J: P&L $
$mm.Cells.Item(16, 11).Value = 100.4         # K: Capital After
$mm.Cells.Item(16, 16).Value = "early_exit"  # P: Exit Reason
$mm.Cells.Item(16, 17).Value = 0.15          # Q: Duration (min)

# New trade #81 appended as sheet row 49.
$mm.Cells.Item(49, 1).Value = 81             # A: Trade #
Set-TextValue $mm.Cells.Item(49, 2) "2026-02-17"  # B: Date
Set-TextValue $mm.Cells.Item(49, 3) "20:53:16"    # C: Time
$mm.Cells.Item(49, 4).Value = "MarketMaking" # D: Strategy
$mm.Cells.Item(49, 5).Value = "UP"           # E: Side
$mm.Cells.Item(49, 6).Value = 0.2            # F: Entry Price
$mm.Cells.Item(49, 8).Value = "OPEN"         # H: Status
$mm.Cells.Item(49, 9).Value = 0              # I: P&L %
$mm.Cells.Item(49, 10).Value = 0             # J: P&L $
$mm.Cells.Item(49, 11).Value = 100.4269627845085 # K: Capital After
$mm.Cells.Item(49, 12).Value = 0             # L: Entry Slippage (bps)
$mm.Cells.Item(49, 13).Value = 0             # M: Exit Slippage (bps)
$mm.Cells.Item(49, 14).Value = 0.6           # N: Confidence
$mm.Cells.Item(49, 15).Value = "Normal spread capture: 19600 bps" # O: Entry Reason
$mm.Cells.Item(49, 17).Value = 0             # Q: Duration (min)
